# Add Sprint Shootout details
# - Several "FP2" session rows are actually "Sprint Shootout" sessions; fix the
#   session label and the corresponding start/end times for those events.
# - Also correct the Qatar GP weekend (FP1/Qualifying/Sprint Race/Race) times.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("F1")
$ws.Activate()

# --- Austria: Sprint Shootout (row 62) ---
$ws.Range("B62").Value = "Sprint Shootout"
$ws.Range("E62").Value = 0.5
$ws.Range("G62").Formula = "=E62+(3*(1/96))"

# --- Belgium: Sprint Shootout (row 77) ---
$ws.Range("B77").Value = "Sprint Shootout"
$ws.Range("E77").Value = 0.5
$ws.Range("G77").Formula = "=E77+(3*(1/96))"

# --- Qatar GP weekend (rows 100-104) ---
# FP1
$ws.Range("E100").Value = 0.6875

# Qualifying
$ws.Range("E101").Value = 0.83333333333333337

# Sprint Shootout
$ws.Range("B102").Value = "Sprint Shootout"

# Sprint Race
$ws.Range("E103").Value = 0.66666666666666663
$ws.Range("G103").Formula = "=E103+(3*(1/96))"

# Race
$ws.Range("E104").Value = 0.83333333333333337

# --- USA: Sprint Shootout (row 107) ---
$ws.Range("B107").Value = "Sprint Shootout"
$ws.Range("E107").Value = 0.52083333333333337
$ws.Range("G107").Formula = "=E107+(3*(1/96))"

# --- Brazil: Sprint Shootout (row 117) ---
$ws.Range("B117").Value = "Sprint Shootout"
$ws.Range("E117").Value = 0.45833333333333331
$ws.Range("G117").Formula = "=E117+(3*(1/96))"

# Leave the selection where the author's last edit was
$null = $ws.Range("G118").Select()
